$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "Datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 23:09"

# --- Swap the order of Fiyi / Dominica rows (205=Dominica, 206=Fiyi before edit) ---
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

# --- Update numeric stats for the affected countries ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2972430
$ws.Range("C4").Value = 36660
$ws.Range("D4").Value = 1285133
$ws.Range("E4").Value = 1554766
$ws.Range("G4").Value = 213
$ws.Range("H4").Value = 132531

# Row 8: Peru
$ws.Range("B8").Value = 302718
$ws.Range("C8").Value = 3638
$ws.Range("D8").Value = 193957
$ws.Range("E8").Value = 98172
$ws.Range("G8").Value = 177
$ws.Range("H8").Value = 10589

# Row 23: Canada
$ws.Range("B23").Value = 105535
$ws.Range("C23").Value = 218
$ws.Range("E23").Value = 27612

# Row 27: Egipto
$ws.Range("B27").Value = 75253
$ws.Range("C27").Value = 1218
$ws.Range("D27").Value = 20726
$ws.Range("E27").Value = 51184
$ws.Range("G27").Value = 63
$ws.Range("H27").Value = 3343

# Row 48: Suiza
$ws.Range("D48").Value = 29300
$ws.Range("E48").Value = 1003

# Row 49: Israel
$ws.Range("B49").Value = 29958
$ws.Range("C49").Value = 788
$ws.Range("D49").Value = 17950
$ws.Range("E49").Value = 11677

# Row 71: Uzbekistan
$ws.Range("B71").Value = 10020
$ws.Range("C71").Value = 312
$ws.Range("D71").Value = 6584
$ws.Range("E71").Value = 3402
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 34

# Row 109: Paraguay
$ws.Range("B109").Value = 2427
$ws.Range("C109").Value = 42
$ws.Range("D109").Value = 1166
$ws.Range("E109").Value = 1241

# Row 131: Jordania
$ws.Range("D131").Value = 942
$ws.Range("E131").Value = 212

# Row 154: Surinam
$ws.Range("B154").Value = 594
$ws.Range("C154").Value = 29
$ws.Range("D154").Value = 292
$ws.Range("E154").Value = 288

# Row 165: Comoras
$ws.Range("B165").Value = 311
$ws.Range("C165").Value = 2
$ws.Range("D165").Value = 266
$ws.Range("E165").Value = 38

# Row 183: Barbados
$ws.Range("B183").Value = 98
$ws.Range("C183").Value = 1
$ws.Range("E183").Value = 1
